$d = $word.ActiveDocument

# 1. "Schedule of testing ... sprints (2-week estimate projection per sprint)..."
#    -> "... sprints (2-3 week estimate projection per sprint)..."
$d.Content.Find.Execute("(2-week estimate", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(2-3 week estimate", 2)

# 2. "...to download each build to playtest Infinity Starship."
#    -> "...to download each build to playtest it."
$d.Content.Find.Execute("to playtest Infinity Starship.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "to playtest it.", 2)

# 3. "Requires satisfaction of windows screens functioning and transitioning correctly from one to another."
#    -> "...from one to another (PlayScreenForm to PlayerIDEntryForm)."
$d.Content.Find.Execute("functioning and transitioning correctly from one to another.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "functioning and transitioning correctly from one to another (PlayScreenForm to PlayerIDEntryForm).", 2)

# 4. Remove the "Player is in main menu (MainMenuWindowForm), where they want to play a new
#    game, clicking play button." bullet entirely (its content is superseded by the next bullet).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Player is in main menu (MainMenuWindowForm)*") {
        $p.Range.Delete()
        break
    }
}

# 5. Rewrite the "When Player gets hit 3 times..." bullet with the new game-over / leaderboard flow.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "When Player gets hit 3 times*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "When Player gets hit 3 times, show message box telling user game’s over with final score.  Transition Player to a small PlayaerIDEntryForm pop up that prompts user to input a valid PlayerID identifier that they can save and store their score in a small leaderboard (Within next game iteration in the PlayerIDEntryForm or in message box pop up) for them to see."
        break
    }
}

# 6. "Player either wants to simply quit by exiting game or play new game."
#    -> "...or play again new game."
$d.Content.Find.Execute("play new game.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "play again new game.", 2)
